$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38, shifting existing rows 38..161 down to 39..162.
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with this week's data point.
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = 44487
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 100112003
$ws.Range("G38").Value = "Ajo"
$ws.Range("H38").Value = "Chino"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 21000
$ws.Range("L38").Value = 21000
$ws.Range("M38").Value = 21000
$ws.Range("N38").Value = "$/caja 10 kilos"
$ws.Range("O38").Value = "China"
$ws.Range("P38").Value = 2100
$ws.Range("Q38").Value = 10
$ws.Range("R38").Value = "Hortaliza"
